# Versión 1.0.5 en PRD - Actualización versión en GIT
#
# "Etapa 4 ES" (sheet 1) had a stray placeholder e-mail
# (icorral@rpatechnologies.es) hard-coded into column B (with a grey-fill
# check-cell style), while the real per-row e-mail addresses were sitting
# unused off in column I. This moves the real e-mail addresses from column
# I into column B (dropping the placeholder + its fill), then removes the
# now-empty column I entirely.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "Etapa 4 ES"
$ws2 = $wb.Worksheets.Item(2)   # "Etapa 4 PT"

# Move the real e-mail addresses (currently parked in column I) into
# column B, replacing the placeholder value, and drop the grey
# "check-cell" fill that was flagging them.
for ($r = 2; $r -le 10; $r++) {
    $email = $ws1.Cells.Item($r, 9).Value2
    $ws1.Cells.Item($r, 2).Value = $email
    $ws1.Cells.Item($r, 2).Interior.Pattern = -4142   # xlPatternNone
}

# Column I is no longer needed now that its data lives in column B.
[void]$ws1.Range("I2:I10").ClearContents()

# Update the saved selections to match where the users left off.
# (Sheet2 is touched first so that sheet1 ends up as the activated /
# tab-selected sheet, matching the saved workbook state.)
$ws2.Activate()
[void]$ws2.Range("J1:J4").Select()

$ws1.Activate()
[void]$ws1.Range("E9").Select()
